$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Workbook-level metadata tweaks (best-effort; many low-level IDs are environment-generated) ---

# --- Populate new address rows 49-90 (styled like the pasted "Normal 2" style) ---
$ws.Cells.Item(49, 1).Value = "Pliar SN  e/  Julio de Cardenas y Lindero, Arroyo Naranjo, La Habana"
$ws.Cells.Item(50, 1).Value = "31 A e/ 320 y 322, La Lisa, La Habana"
$ws.Cells.Item(51, 1).Value = "C  e/  Delicias y San Francisco, San  Miguel del Padron, La Habana"
$ws.Cells.Item(52, 1).Value = "39 e/ San Juan Bautista y Union, Plaza de la Revolucion, La Habana"
$ws.Cells.Item(53, 1).Value = "Línea e/ 20 y 22, Vedado, Plaza de la Revolucion, La Habana"
$ws.Cells.Item(54, 1).Value = "1ra e/ D y E, reparto Luyano Moderno, San  Miguel del Padron, La Habana"
$ws.Cells.Item(55, 1).Value = "29 e/ 310 y 312, La Lisa, La Habana"
$ws.Cells.Item(56, 1).Value = "428 e/ 19 y 21, Pena altas, Guanabo, La Habana del Este, La Habana"
$ws.Cells.Item(57, 1).Value = "246 e/ 33 C y 35, La Lisa, La Habana"
$ws.Cells.Item(58, 1).Value = "Figura e/ Vives y Esperanza, La Habana Vieja, La Habana"
$ws.Cells.Item(59, 1).Value = "B e/ 1 y 2, Mercedita, San Miguel del Padron, La Habana"
$ws.Cells.Item(60, 1).Value = "161 e/ 310 y 314, Valle Brande, La Lisa, La Habana"
$ws.Cells.Item(61, 1).Value = "Santa Ana e/ Luco y Villanueva, Diez de Octubre, La Habana"
$ws.Cells.Item(62, 1).Value = "44 e/ 19 Y  21, Playa, La Habana"
$ws.Cells.Item(63, 1).Value = " Joaquin Delgado e Santa Clara y Esperanza, reparto Parraga, Arroyo Naranjo, La Habana"
$ws.Cells.Item(64, 1).Value = "CALLE C  # 9512 A  /  6 y 10 Reparto Altahabana ,BOYEROS,La Habana"
$ws.Cells.Item(65, 1).Value = "Calle 160 No 4509E / Ave 45 y Ave 47,LA LISA,La Habana"
$ws.Cells.Item(66, 1).Value = "Calle Maceo  # 61 Alto entre Bertenati y Nazareno,PLAZA DE LA REVOLUCION,La Habana"
$ws.Cells.Item(67, 1).Value = "Destramoes  # 58 entre Luis esteves y Lacret, Santo Suarez,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(68, 1).Value = "Calle Máximo Gómez  # 264 Altos  E / 27 DE Noviembre y Pereira,REGLA,La Habana"
$ws.Cells.Item(69, 1).Value = "Calle San Mariano  # 761 entre San Juan Bosco y Graciela apto 8, Vibora,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(70, 1).Value = "Calle 111 E / cALLE 24 y Calle 24 A Edificio 6 apt 4 Reparto Sierra Maestra,BOYEROS,La Habana"
$ws.Cells.Item(71, 1).Value = "Calle Barbería e /  Calle Universidad y Calle Estevez edf 32 apto 4,CERRO,La Habana"
$ws.Cells.Item(72, 1).Value = "CALLE 18  #  505 BAJO ENTRE CONCEPCION Y SAN FRANCISCO LAWTON,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(73, 1).Value = "180  #  42109 ENTRE 421 Y 423, LA AURORA,BOYEROS,La Habana"
$ws.Cells.Item(74, 1).Value = "Calle Gertrudis Oeste  # 459 E /  Anita y Goicuria Repto Sevillano,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(75, 1).Value = "Avenida Santa Catalina e / Mayia Rodriguez  y La Sola,Edificio 616,Apto 3,Reparto Santos Suarez,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(76, 1).Value = "Calle Pinar  # 8003 E / Santa Teresa y Collazo,Reparto Ponce,ARROYO NARANJO,La Habana"
$ws.Cells.Item(77, 1).Value = "Avenida Ciudamar  # 17521 e /  Calle 13 y Calle 15, Reparto Ciudamar,SAN MIGUEL DEL PADRON,La Habana"
$ws.Cells.Item(78, 1).Value = "Calzada de 10 de Octubre  # 1155 apto 1 e /  Santa Catalina y San Mariano,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(79, 1).Value = "Calle 5ta B  #  9606 2 int,entre calle 96 y 98, Barrio Querejeta,PLAYA,La Habana"
$ws.Cells.Item(80, 1).Value = "Calle 1 No  17416 E / Calle A y Calle San Luis,Reparto Encanto,SAN MIGUEL DEL PADRON,La Habana"
$ws.Cells.Item(81, 1).Value = "Calle 9na  # 14 ,Apto 3 E / Calle E y Calle D,Barrio Lawton,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(82, 1).Value = "Calle Tenerifr # 156 apto  6 E /  Rastro y Carmen  ,CENTRO HABANA,La Habana"
$ws.Cells.Item(83, 1).Value = "Ave Entrada E / Ave de los Ocujes y Ave  de la Ceiba Edificio 2 APTO 5 Reparto Santa Catalina,CERRO,La Habana"
$ws.Cells.Item(84, 1).Value = "Sta Catalina e / Vento y San Juán Bosco ed 817 apto 4,DIEZ DE OCTUBRE,La Habana"
$ws.Cells.Item(85, 1).Value = "Aliados e/ Pasaje D y Central, San Miguel del Padrón, La Habana"
$ws.Cells.Item(86, 1).Value = "Estancia e/ San Pedro y Lombill, Nuevo Vedado, Plaza de la Revolucion, La Habana"
$ws.Cells.Item(87, 1).Value = "Estrella e/ Aguila y Angeles, Centro Habana, La Habana"
$ws.Cells.Item(88, 1).Value = "Virginia e/ Pinar del Rio y Woodbury, reparto Callejas, Arroyo Naranjo, La Habana"
$ws.Cells.Item(89, 1).Value = "Recurso e/ Masarredo y Lindero, Nuevo Vedado, Plaza de la Revolucion, La Habana"
$ws.Cells.Item(90, 1).Value = "San Quintin e/ Salvador y Cerezo, Cerro, La Habana"

# --- Populate new address rows 91-100 (default style, no special font) ---
$ws.Cells.Item(91, 1).Value = "Calle 6ta entre 16 y 17 Edificio 37 Apto 29 Reparto Guiteras  LA HABANA DEL ESTE La Habana"
$ws.Cells.Item(92, 1).Value = "Calle 308   1904A entre 19 y 21 Reparto Santa Fe  PLAYA La Habana"
$ws.Cells.Item(93, 1).Value = "Calle Santa Rosa     8 e entre San Antonio y Rizo  PLaza PLAZA DE LA REVOLUCION La Habana"
$ws.Cells.Item(94, 1).Value = "Calle 17    853 e entre  4 y 6 PLAZA DE LA REVOLUCION La Habana"
$ws.Cells.Item(95, 1).Value = "67A No 11404 e 114 y 116 MARIANAO La Habana"
$ws.Cells.Item(96, 1).Value = "Calle 203 entre 290 y Prensa Latina, Poblado Pueblo del Chico, BOYEROS, La Habana"
$ws.Cells.Item(97, 1).Value = "A # 319/11 Y 12"
$ws.Cells.Item(98, 1).Value = "REAL # 181"
$ws.Cells.Item(99, 1).Value = "JACINTO ROY # 14 / ANGELITA Y BERENGUER, VIEJA LINDA"
$ws.Cells.Item(100, 1).Value = "PUERTA CERRADA 222 / ALAMBIQUE Y FLORIDA"

# --- Apply the "pasted" style (Normal 2 cell style + black Calibri font) to rows 49-90 ---
$newStyle = $wb.Styles.Add("Normal 2")
$newStyle.Font.Color = 0
$ws.Range("A49:A90").Style = "Normal 2"

# --- Update selection / view to mirror the final sheet view state ---
$ws.Range("A97:A100").Select()

Write-Host "Applied 100-address update"
